$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "INDICATIVE FOREIGN EXCHANGE RATES AS AT 07-Feb-2025"

# Shift "ND: NO DEALINGS" / "DISCLAIMER" block down by 23 rows to make room for new data
$ws.Range("A11:K13").Cut($ws.Range("A34:K36"))

# Row 10: AUSTRALIA
$ws.Cells.Item(10, 1).Value = "AUSTRALIA"
$ws.Cells.Item(10, 2).Value = "Australia Dollars"
$ws.Cells.Item(10, 3).Value = "AUD"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 29.36
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 29.01
$ws.Cells.Item(10, 8).Value = 30.8
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 30.8
$ws.Cells.Item(10, 11).Value = "07-Feb-2025 08:56"

# Row 11: CANADA
$ws.Cells.Item(11, 1).Value = "CANADA"
$ws.Cells.Item(11, 2).Value = "Canadian Dollar"
$ws.Cells.Item(11, 3).Value = "CAD"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 32.75
$ws.Cells.Item(11, 6).Value = 32.65
$ws.Cells.Item(11, 7).Value = 32.26
$ws.Cells.Item(11, 8).Value = 34.35
$ws.Cells.Item(11, 9).Value = 34.35
$ws.Cells.Item(11, 10).Value = 34.35
$ws.Cells.Item(11, 11).Value = "07-Feb-2025 08:56"

# Row 12: CHINA
$ws.Cells.Item(12, 1).Value = "CHINA"
$ws.Cells.Item(12, 2).Value = "China Yuan Renminbi"
$ws.Cells.Item(12, 3).Value = "CNY"
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 6.36
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 6.69
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = "07-Feb-2025 08:56"

# Row 13: EUR COUNTRIES
$ws.Cells.Item(13, 1).Value = "EUR COUNTRIES"
$ws.Cells.Item(13, 2).Value = "Euro"
$ws.Cells.Item(13, 3).Value = "EUR"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 47.83
$ws.Cells.Item(13, 6).Value = 47.71
$ws.Cells.Item(13, 7).Value = 47.37
$ws.Cells.Item(13, 8).Value = 49.26
$ws.Cells.Item(13, 9).Value = 49.26
$ws.Cells.Item(13, 10).Value = 49.26
$ws.Cells.Item(13, 11).Value = "07-Feb-2025 08:56"

# Row 14: GREAT BRITAIN
$ws.Cells.Item(14, 1).Value = "GREAT BRITAIN"
$ws.Cells.Item(14, 2).Value = "Great Britain Pound"
$ws.Cells.Item(14, 3).Value = "GBP"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 57.24
$ws.Cells.Item(14, 6).Value = 57.1
$ws.Cells.Item(14, 7).Value = 56.683
$ws.Cells.Item(14, 8).Value = 58.95
$ws.Cells.Item(14, 9).Value = 58.95
$ws.Cells.Item(14, 10).Value = 58.95
$ws.Cells.Item(14, 11).Value = "07-Feb-2025 08:56"

# Row 15: HONG KONG
$ws.Cells.Item(15, 1).Value = "HONG KONG"
$ws.Cells.Item(15, 2).Value = "Hong Kong Dollar"
$ws.Cells.Item(15, 3).Value = "HKD"
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(15, 5).Value = 5.89
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 5.73
$ws.Cells.Item(15, 8).Value = 6.23
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 6.23
$ws.Cells.Item(15, 11).Value = "07-Feb-2025 08:56"

# Row 16: INDIA
$ws.Cells.Item(16, 1).Value = "INDIA"
$ws.Cells.Item(16, 2).Value = "India Rupee"
$ws.Cells.Item(16, 3).Value = "INR"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0.5347
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0.5689
$ws.Cells.Item(16, 9).Value = 0.5689
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = "07-Feb-2025 08:56"

# Row 17: JAPAN
$ws.Cells.Item(17, 1).Value = "JAPAN"
$ws.Cells.Item(17, 2).Value = "Japan Yen"
$ws.Cells.Item(17, 3).Value = "JPY"
$ws.Cells.Item(17, 4).Value = 100
$ws.Cells.Item(17, 5).Value = 30.81
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 30.04
$ws.Cells.Item(17, 8).Value = 32.44
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 32.44
$ws.Cells.Item(17, 11).Value = "07-Feb-2025 08:56"

# Row 18: KENYA
$ws.Cells.Item(18, 1).Value = "KENYA"
$ws.Cells.Item(18, 2).Value = "Kenya Shilling"
$ws.Cells.Item(18, 3).Value = "KES"
$ws.Cells.Item(18, 4).Value = 100
$ws.Cells.Item(18, 5).Value = 36.57
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 38.11
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = "07-Feb-2025 08:56"

# Row 19: MADAGASCAR
$ws.Cells.Item(19, 1).Value = "MADAGASCAR"
$ws.Cells.Item(19, 2).Value = "Malagasy Ariary"
$ws.Cells.Item(19, 3).Value = "MGA"
$ws.Cells.Item(19, 4).Value = 100
$ws.Cells.Item(19, 5).Value = 1.01
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 1.07
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = "07-Feb-2025 08:56"

# Row 20: NEW ZEALAND
$ws.Cells.Item(20, 1).Value = "NEW ZEALAND"
$ws.Cells.Item(20, 2).Value = "Cook Islands New Zealand Dollars"
$ws.Cells.Item(20, 3).Value = "NZD"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 26.16
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 25.7
$ws.Cells.Item(20, 8).Value = 27.39
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 27.39
$ws.Cells.Item(20, 11).Value = "07-Feb-2025 08:56"

# Row 21: NORWAY
$ws.Cells.Item(21, 1).Value = "NORWAY"
$ws.Cells.Item(21, 2).Value = "Bouvet Island Norway Kroner"
$ws.Cells.Item(21, 3).Value = "NOK"
$ws.Cells.Item(21, 4).Value = 1
$ws.Cells.Item(21, 5).Value = 4.05
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 3.93
$ws.Cells.Item(21, 8).Value = 4.28
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 4.28
$ws.Cells.Item(21, 11).Value = "07-Feb-2025 08:56"

# Row 22: PAKISTAN
$ws.Cells.Item(22, 1).Value = "PAKISTAN"
$ws.Cells.Item(22, 2).Value = "Pakistan Rupees"
$ws.Cells.Item(22, 3).Value = "PKR"
$ws.Cells.Item(22, 4).Value = 100
$ws.Cells.Item(22, 5).Value = 16.5
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 17.47
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = "07-Feb-2025 08:56"

# Row 23: SAUDI ARABIA
$ws.Cells.Item(23, 1).Value = "SAUDI ARABIA"
$ws.Cells.Item(23, 2).Value = "Saudi Arabia Riyals"
$ws.Cells.Item(23, 3).Value = "SAR"
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(23, 5).Value = 12.28
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 11.97
$ws.Cells.Item(23, 8).Value = 12.99
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 12.99
$ws.Cells.Item(23, 11).Value = "07-Feb-2025 08:56"

# Row 24: SEYCHELLES
$ws.Cells.Item(24, 1).Value = "SEYCHELLES"
$ws.Cells.Item(24, 2).Value = "Seychelles Rupees"
$ws.Cells.Item(24, 3).Value = "SCR"
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 3.17
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 2.92
$ws.Cells.Item(24, 8).Value = 3.47
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 3.47
$ws.Cells.Item(24, 11).Value = "07-Feb-2025 08:56"

# Row 25: SINGAPORE
$ws.Cells.Item(25, 1).Value = "SINGAPORE"
$ws.Cells.Item(25, 2).Value = "Singapore Dollars"
$ws.Cells.Item(25, 3).Value = "SGD"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 33.84
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 33.59
$ws.Cells.Item(25, 8).Value = 35.63
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 35.63
$ws.Cells.Item(25, 11).Value = "07-Feb-2025 08:56"

# Row 26: SOUTH AFRICA
$ws.Cells.Item(26, 1).Value = "SOUTH AFRICA"
$ws.Cells.Item(26, 2).Value = "South Africa Rand"
$ws.Cells.Item(26, 3).Value = "ZAR"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 2.52
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 2.47
$ws.Cells.Item(26, 8).Value = 2.68
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 2.68
$ws.Cells.Item(26, 11).Value = "07-Feb-2025 08:56"

# Row 27: SRI LANKA
$ws.Cells.Item(27, 1).Value = "SRI LANKA"
$ws.Cells.Item(27, 2).Value = "Sri Lanka Rupee"
$ws.Cells.Item(27, 3).Value = "LKR"
$ws.Cells.Item(27, 4).Value = 100
$ws.Cells.Item(27, 5).Value = 15.84
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 16.58
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = "07-Feb-2025 08:56"

# Row 28: SWEDEN
$ws.Cells.Item(28, 1).Value = "SWEDEN"
$ws.Cells.Item(28, 2).Value = "Sweden Kronor"
$ws.Cells.Item(28, 3).Value = "SEK"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 5).Value = 4.18
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 4.06
$ws.Cells.Item(28, 8).Value = 4.38
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 4.38
$ws.Cells.Item(28, 11).Value = "07-Feb-2025 08:56"

# Row 29: SWITZERLAND
$ws.Cells.Item(29, 1).Value = "SWITZERLAND"
$ws.Cells.Item(29, 2).Value = "Switzerland Franc"
$ws.Cells.Item(29, 3).Value = "CHF"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = 51.09
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 50.58
$ws.Cells.Item(29, 8).Value = 53.8
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 53.8
$ws.Cells.Item(29, 11).Value = "07-Feb-2025 08:56"

# Row 30: TANZANIA
$ws.Cells.Item(30, 1).Value = "TANZANIA"
$ws.Cells.Item(30, 2).Value = "Tanzania Shillings"
$ws.Cells.Item(30, 3).Value = "TZS"
$ws.Cells.Item(30, 4).Value = 100
$ws.Cells.Item(30, 5).Value = 1.81
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 1.9
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = "07-Feb-2025 08:56"

# Row 31: UNITED ARAB EMIRATES
$ws.Cells.Item(31, 1).Value = "UNITED ARAB EMIRATES"
$ws.Cells.Item(31, 2).Value = "UAE Dirham"
$ws.Cells.Item(31, 3).Value = "AED"
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(31, 5).Value = 12.86
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 11.91
$ws.Cells.Item(31, 8).Value = 13.46
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 13.46
$ws.Cells.Item(31, 11).Value = "07-Feb-2025 08:56"

# Row 32: USA
$ws.Cells.Item(32, 1).Value = "USA"
$ws.Cells.Item(32, 2).Value = "United States Dollars"
$ws.Cells.Item(32, 3).Value = "USD"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = 46.29
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 46.09
$ws.Cells.Item(32, 8).Value = 47.19
$ws.Cells.Item(32, 9).Value = 47.19
$ws.Cells.Item(32, 10).Value = 47.662
$ws.Cells.Item(32, 11).Value = "07-Feb-2025 08:56"
